$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.786.77"
$ws.Range("E2").Value = "  +2.77%  "
$ws.Range("D3").Value = "3.039.06"
$ws.Range("E3").Value = "  +2.36%  "
$ws.Range("E4").Value = "  +0.16%  "
$c = $ws.Range("D5")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "515.14"
$c.Style = $origStyle
$ws.Range("E5").Value = "  +2.71%  "
$c = $ws.Range("D6")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "140.03"
$c.Style = $origStyle
$ws.Range("E6").Value = "  +4.30%  "
$c = $ws.Range("D7")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = $origStyle
$ws.Range("E7").Value = "  +0.02%  "
$c = $ws.Range("D8")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.442"
$c.Style = $origStyle
$ws.Range("E8").Value = "  +3.35%  "
$c = $ws.Range("D9")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "7.48"
$c.Style = $origStyle
$ws.Range("E9").Value = "  +2.26%  "
$ws.Range("E10").Value = "  +4.07%  "
$c = $ws.Range("D11")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.367"
$c.Style = $origStyle
$ws.Range("E11").Value = "  +5.04%  "
$ws.Range("D12").Value = "3.571.28"
$ws.Range("E12").Value = "  +2.64%  "
$ws.Range("E13").Value = "  +2.35%  "
$c = $ws.Range("D14")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "26.67"
$c.Style = $origStyle
$ws.Range("E14").Value = "  +5.91%  "
$c = $ws.Range("D15")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.0000167"
$c.Style = $origStyle
$ws.Range("E15").Value = "  +10.52%  "
$ws.Range("D16").Value = "57.832.64"
$ws.Range("E16").Value = "  +2.81%  "
$c = $ws.Range("D17")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "6.20"
$c.Style = $origStyle
$ws.Range("E17").Value = "  +9.76%  "
$ws.Range("D18").Value = "3.034.80"
$ws.Range("E18").Value = "  +2.15%  "
$c = $ws.Range("D19")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "12.95"
$c.Style = $origStyle
$ws.Range("E19").Value = "  +5.46%  "
$c = $ws.Range("D20")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "8.01"
$c.Style = $origStyle
$ws.Range("E20").Value = "  +3.61%  "
$c = $ws.Range("D21")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "335.05"
$c.Style = $origStyle
$ws.Range("E21").Value = "  +4.17%  "
$c = $ws.Range("D22")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "5.78"
$c.Style = $origStyle
$ws.Range("E22").Value = "  +1.65%  "
$c = $ws.Range("D23")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.998"
$c.Style = $origStyle
$ws.Range("E23").Value = "  +0.00%  "
$c = $ws.Range("D24")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.499"
$c.Style = $origStyle
$ws.Range("E24").Value = "  +6.45%  "
$c = $ws.Range("D25")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "64.86"
$c.Style = $origStyle
$ws.Range("E25").Value = "  +4.85%  "
$c = $ws.Range("D26")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.168"
$c.Style = $origStyle
$ws.Range("E26").Value = "  +3.87%  "
$ws.Range("B27").Value = "Binance-PegBSC-USD"
$ws.Range("C27").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$c = $ws.Range("D27")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = $origStyle
$ws.Range("E27").Value = "  +0.14%  "
$ws.Range("B28").Value = "PEPE"
$ws.Range("C28").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D28").Value = "0.0₃0937"
$ws.Range("E28").Value = "  +6.06%  "
$c = $ws.Range("D29")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "6.82"
$c.Style = $origStyle
$ws.Range("E29").Value = "  +5.34%  "
$c = $ws.Range("D30")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "7.45"
$c.Style = $origStyle
$ws.Range("E30").Value = "  +10.28%  "
$c = $ws.Range("D31")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "1.81"
$c.Style = $origStyle
$ws.Range("E31").Value = "  +4.12%  "
$c = $ws.Range("D32")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "1.22"
$c.Style = $origStyle
$ws.Range("E32").Value = "  +2.84%  "
$c = $ws.Range("D33")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "20.89"
$c.Style = $origStyle
$ws.Range("E33").Value = "  +2.58%  "
$c = $ws.Range("D34")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "155.98"
$c.Style = $origStyle
$ws.Range("E34").Value = "  -1.44%  "
$c = $ws.Range("D35")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "4.72"
$c.Style = $origStyle
$ws.Range("E35").Value = "  +6.19%  "
$c = $ws.Range("D36")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "5.86"
$c.Style = $origStyle
$ws.Range("E36").Value = "  +6.30%  "
$c = $ws.Range("D37")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "1.28"
$c.Style = $origStyle
$ws.Range("E37").Value = "  +2.21%  "
$c = $ws.Range("D38")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "24.86"
$c.Style = $origStyle
$ws.Range("E38").Value = "  +7.94%  "
$c = $ws.Range("D39")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.0687"
$c.Style = $origStyle
$ws.Range("E39").Value = "  +2.23%  "
$ws.Range("D40").Value = "3.079.67"
$ws.Range("E40").Value = "  +2.51%  "
$c = $ws.Range("D41")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "37.61"
$c.Style = $origStyle
$ws.Range("E41").Value = "  +4.03%  "
$c = $ws.Range("D42")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "3.87"
$c.Style = $origStyle
$ws.Range("E42").Value = "  +9.48%  "
$ws.Range("E43").Value = "  +0.18%  "
$c = $ws.Range("D44")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.660"
$c.Style = $origStyle
$ws.Range("E44").Value = "  +3.22%  "
$ws.Range("D45").Value = "2.304.44"
$ws.Range("E45").Value = "  +2.71%  "
$c = $ws.Range("D46")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "1.43"
$c.Style = $origStyle
$ws.Range("E46").Value = "  +3.20%  "
$c = $ws.Range("D47")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.992"
$c.Style = $origStyle
$ws.Range("E47").Value = "  +1.30%  "
$c = $ws.Range("D48")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "6.04"
$c.Style = $origStyle
$ws.Range("E48").Value = "  +5.25%  "
$c = $ws.Range("D49")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.0240"
$c.Style = $origStyle
$ws.Range("E49").Value = "  +3.01%  "
$c = $ws.Range("D50")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "19.54"
$c.Style = $origStyle
$ws.Range("E50").Value = "  +3.27%  "
$c = $ws.Range("D51")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "1.82"
$c.Style = $origStyle
$ws.Range("E51").Value = "  -5.30%  "
